# OperationScenario_Component_Boiler.xlsx edit
# The "heating_element_power" and "heating_element_power_unit" columns
# (E and F) are removed from the boiler component table; the remaining
# columns (carnot_efficiency_factor, heating_supply_temperature,
# hot_water_supply_temperature) shift left to become E, F, G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns E:F (heating_element_power, heating_element_power_unit).
# This shifts the following columns (carnot_efficiency_factor,
# heating_supply_temperature, hot_water_supply_temperature) left.
$ws.Range("E1:F1").EntireColumn.Delete()

# Resize the now-shifted columns to fit their (longer) header text.
$ws.Columns.Item(5).AutoFit()
$ws.Columns.Item(6).AutoFit()
$ws.Columns.Item(7).AutoFit()

# Match the saved selection/active cell of the sheet.
$ws.Range("E7").Select()
